$wb = $excel.ActiveWorkbook

# --- Login Page (sheet1): update Corporate Email and Admin Email values ---
$wsLogin = $wb.Worksheets.Item("Login Page")
$wsLogin.Range("B2").Value = "wipro@mailinator.com"
$wsLogin.Range("B4").Value = "adminsiva@nada.email"
$wsLogin.Range("B5").Select()

# --- Write data (sheet5): update Job ID / CorpDashGS / CorpDashGO values ---
$wsWrite = $wb.Worksheets.Item("Write data")
$wsWrite.Range("B2").Value = "1238"
$wsWrite.Range("B3").Value = "Active JDs - 94; Profiles Received - 318; Expenses - 1528750"
$wsWrite.Range("B4").Value = "Offered - 28; Onboarded - 26; TotalMoneySaved - 0"
